$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / already-non-numeric string cells: direct assignment ---
$ws.Range("D2").Value = "30.208.90"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.861.65"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  +11.76%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "1.865.99"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "30.202.41"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +7.57%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.108.82"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("E35").Value = "  +3.68%  "
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("E39").Value = "  +4.53%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("E41").Value = "  +3.74%  "
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("E51").Value = "  +0.27%  "

# --- Numeric-looking strings in the Price column: force text storage so
#     Excel keeps them as strings (matching source data) instead of casting
#     to Double, while keeping the cell style index untouched. ---
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "236.23"
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4677"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2848"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06526"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "21.78"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07898"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "97.34"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.161"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.6797"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "279.05"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "13.50"
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.000007311"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.359"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.9999"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.160"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "168.20"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.234"
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.932"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.380"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.09738"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.479"
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.049"
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.04725"
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.135"
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7072"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.706"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01865"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.623"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "6.305"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.952"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.8484"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.4175"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "103.39"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "969.01"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.189"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.269"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "34.10"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.05640"
$cell.Style = "Normal"
